$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert three new worksheets between "Configs" and
#    "Possible policies combinations": Landuse Dict, Policy Dict, Combination
# ---------------------------------------------------------------------------
$configs = $wb.Worksheets.Item("Configs")

$landuse = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $configs)
$landuse.Name = "Landuse Dict"

$policyDict = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $landuse)
$policyDict.Name = "Policy Dict"

$combination = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $policyDict)
$combination.Name = "Combination"

# ---------------------------------------------------------------------------
# 2. Fill "Landuse Dict" - landuse category -> numeric code
# ---------------------------------------------------------------------------
$landuse.Range("A1").Value = "N.A."
$landuse.Range("B1").Value = -1
$landuse.Range("A2").Value = "Residential"
$landuse.Range("B2").Value = 0
$landuse.Range("A3").Value = "Commercial"
$landuse.Range("B3").Value = 1
$landuse.Range("A4").Value = "Mixed-use"
$landuse.Range("B4").Value = 2
$landuse.Range("A5").Value = "Green"
$landuse.Range("B5").Value = 3
$landuse.Columns.Item(1).ColumnWidth = 10

# ---------------------------------------------------------------------------
# 3. Fill "Policy Dict" - policy name -> numeric code
# ---------------------------------------------------------------------------
$policyDict.Range("A1").Value = "Urban gardening"
$policyDict.Range("B1").Value = 1
$policyDict.Range("A2").Value = "Limited land allocation for fodder crops"
$policyDict.Range("B2").Value = 2
$policyDict.Range("A3").Value = "Sustainable farming production system"
$policyDict.Range("B3").Value = 3
$policyDict.Range("A4").Value = "Draining garden design"
$policyDict.Range("B4").Value = 4
$policyDict.Range("A5").Value = "Rainwater harvesting"
$policyDict.Range("B5").Value = 5
$policyDict.Range("A6").Value = "On-site wastewater purification"
$policyDict.Range("B6").Value = 6
$policyDict.Range("A7").Value = "Solar power roofs"
$policyDict.Range("B7").Value = 7
$policyDict.Range("A8").Value = "Energy-saving households behavior"
$policyDict.Range("B8").Value = 8
$policyDict.Range("A9").Value = "Biomass efficiency improvement"
$policyDict.Range("B9").Value = 9
$policyDict.Range("A10").Value = "Wind power"
$policyDict.Range("B10").Value = 10
$policyDict.Columns.Item(1).ColumnWidth = 33

# ---------------------------------------------------------------------------
# 4. Fill "Combination" - numeric code -> combination of policy numbers
# ---------------------------------------------------------------------------
$combination.Range("A1").Value = 12
$combination.Range("B1").Value = "1,3"
$combination.Range("A2").Value = 13
$combination.Range("B2").Value = "1,4"
$combination.Range("A3").Value = 14
$combination.Range("B3").Value = "2,3"
$combination.Range("A4").Value = 15
$combination.Range("B4").Value = "2,3,9"
$combination.Range("A5").Value = 16
$combination.Range("B5").Value = "5,6"
$combination.Range("A6").Value = 17
$combination.Range("B6").Value = "5,6,7"
$combination.Range("A7").Value = 18
$combination.Range("B7").Value = "7,8"
$combination.Range("A8").Value = 19
$combination.Range("B8").Value = "5,6,7,8"

# ---------------------------------------------------------------------------
# 5. Adjust the "Configs" sheet: widen column B, move its selection
# ---------------------------------------------------------------------------
$configs.Columns.Item(2).ColumnWidth = 12.5

# ---------------------------------------------------------------------------
# 6. Restore per-sheet selections / active-cell state to match the target
#    view. Each .Select() call both sets the sheet's remembered selection
#    and implicitly activates that sheet; the LAST sheet activated below
#    ("Combination") ends up as the workbook's active tab, matching the
#    target activeTab.
# ---------------------------------------------------------------------------
$distance = $wb.Worksheets.Item("Distance")
[void]$distance.Range("I20").Select()

$configs.Activate()
[void]$configs.Range("A2").Select()

[void]$landuse.Range("B6").Select()

[void]$policyDict.Range("A1").Select()

[void]$combination.Range("C5").Select()

Write-Output "Workbook restructuring complete."
